# Revisao na lista de materiais: adiciona item "Perfil T-slot 60x30 compr. 500mm"
# (qtd. 4) logo apos o "Perfil T-slot 60x30 compr. 300mm" na planilha de BOM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha em branco na posicao 17, empurrando o restante da
# lista (Parafusos, Fabricacao, Eletronica, Estrutura, ...) uma linha abaixo.
$ws.Rows(17).Insert()

# Preenche a nova linha com o novo item e sua quantidade.
$ws.Range("B17").Value = "Perfil T-slot 60x30 compr. 500mm"
$ws.Range("C17").Value = 4

# Atualiza a celula selecionada / topo visivel da planilha.
$ws.Range("C18").Select()
